$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing the existing rows 5-12 down to 6-13.
$ws.Rows("5:5").Insert()

# Populate the new row 5 with the new weekly data point.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44614
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101004
$ws.Range("J5").Value = "Frambuesa"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Provincia de Linares"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 2
